# "large code cleanup and addition of some indicators"
#
# The sheet already contains Bollinger-Band style indicator formulas
# (AVERAGE / STDEVPA based) in columns C:E starting at row 20 (the first
# row where a full 20-period lookback window is available). Rows 1-19
# were left blank for those columns; this change backfills them with 0
# so the indicator columns are fully populated down the sheet, and
# restores the active cell/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Backfill the warm-up rows (1-19) of the indicator columns C, D and E with 0.
$ws.Range("C1:E19").Value = 0

# Move the selection/active cell to D20 (first real indicator cell).
$null = $ws.Range("D20").Select()
